$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy H1's formatting (bold, border,
# centered) so they reuse the existing header style rather than minting a new one.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new columns I (I0) and J (IF), rows 2-15
$values = @{
    2  = @(3, 4)
    3  = @(6, 7)
    4  = @(5, 6)
    5  = @(7, 7)
    6  = @(10, 10)
    7  = @(5, 5)
    8  = @(7, 7)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(6, 8)
    14 = @(1, 3)
    15 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
